$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting rows 188:272 down to 189:273
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with its values
$ws.Cells.Item(188, 1).Value = 6
$ws.Cells.Item(188, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44875
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = 100112022
$ws.Cells.Item(188, 7).Value = "Arveja Verde"
$ws.Cells.Item(188, 8).Value = "Perfection"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 320
$ws.Cells.Item(188, 11).Value = 18000
$ws.Cells.Item(188, 12).Value = 20000
$ws.Cells.Item(188, 13).Value = 18938
$ws.Cells.Item(188, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(188, 16).Value = 758
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Apply the date style (s="2") to D188, matching the D column's existing style
$ws.Cells.Item(188, 4).NumberFormat = $ws.Cells.Item(189, 4).NumberFormat
